# Adicion de carrera a tabla usuarios
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: E1 = "Carrera" (reuse header formatting from D1) ---
$ws.Range("E1").Value = "Carrera"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# --- Row 2 (Admin General / JEFE_CARRERA) - no carrera, just formatted like body ---
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# --- Row 3 (Juan Perez / PROFESOR) - no carrera ---
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# --- Row 4 (Maria Silva / PROFESOR) - no carrera ---
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# --- Row 5 (Carlos Diaz / ALUMNO) ---
$ws.Range("E5").Value = "Ingeniería Civil en Informática"
$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4122)

# --- Row 6 (Ana Torres / ALUMNO) ---
$ws.Range("E6").Value = "Ingeniería de Ejecución en Computación"
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)

# --- Row 7 (Fabian Test / ALUMNO) ---
$ws.Range("E7").Value = "Ingeniería Civil Industrial"
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)

# --- Row 8 (Sofia Lagos / ALUMNO) ---
$ws.Range("E8").Value = "Ingeniería Civil en Informática"
$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column E width ---
$ws.Columns.Item(5).ColumnWidth = 35.75

# --- Selection moves to E8, matching last-edited cell ---
$ws.Range("E8").Select()

Write-Output "edit complete"
